# Fix the "Delay less than 0 second" bullet on the hedging cheat sheet slide:
#   "Delay less than 0 second " -> "Delay less than " + "0 seconds " (two runs)
# i.e. pluralise "second" -> "seconds" while keeping the existing run boundary
# pattern used by the rest of the bullet list (a separate run for the numeric
# condition clause).

$p = $ppt.ActivePresentation

$oldPhrase = "Delay less than 0 second "
$firstPart = "Delay less than "
$secondPart = "0 seconds "

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        $fullText = $tr.Text
        $charIdx = $fullText.IndexOf($oldPhrase)
        if ($charIdx -ge 0) {
            # COM TextRange/Characters collections are 1-based.
            $start = $charIdx + 1

            # Re-point the tail of the existing run ("0 second ") to read
            # "0 seconds " - this also splits it off into its own run, just
            # like the first part ("Delay less than ") already is.
            $tailStart = $start + $firstPart.Length
            $tailLength = $oldPhrase.Length - $firstPart.Length
            $tail = $tr.Characters($tailStart, $tailLength)
            $tail.Text = $secondPart
        }
    }
}
